$wb = $excel.ActiveWorkbook

$metadata = $wb.Worksheets.Item("Metadata")
$metadata.Range("B8").Value = "2023-09-15T14:15:22+00:00"

$concepts = $wb.Worksheets.Item("Concepts")
$concepts.Range("C2").Value = "Phase 3 and phase 4"
$concepts.Range("D2").Value = "Trials that are a combination of phases 3 and 4."
